# Case_5_190 (380 kV) results: update res_bus vm_pu.xlsx voltage magnitudes (p.u.)
# Columns B-F and I-M, rows 2-25 (bus voltage results for the 380 kV scenario).
# Columns G, H, N and row 1 (headers) are left unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    "B" = @(1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02)
    "C" = @(1.017506523434188, 1.018562034851111, 1.019244696683664, 1.019531611501204, 1.019579781335644, 1.019248530751122, 1.017863306583473, 1.01541981036163, 1.013788997663064, 1.013082384537353, 1.012819845721687, 1.012876164445788, 1.013060684462165, 1.013174363832497, 1.013835883372152, 1.014250712373064, 1.014492630823719, 1.014575111313388, 1.014206209743933, 1.013006349899533, 1.012251537196672, 1.01265171755005, 1.014226318703305, 1.016051825427843)
    "D" = @(1.02335253314526, 1.024119870061698, 1.024615847283845, 1.024824226797496, 1.024859207058205, 1.02461863216896, 1.023611970419874, 1.02183394875675, 1.020645788469741, 1.020130629594745, 1.019939174276089, 1.01998024673383, 1.020114805930144, 1.020197698709045, 1.020679963542554, 1.020982293713324, 1.021158572445819, 1.021218667871703, 1.020949863310586, 1.020075184431324, 1.019524646500894, 1.019816553440379, 1.020964517415652, 1.022294104261442)
    "E" = @(1.018887464434443, 1.019785563326101, 1.020366789292499, 1.020611160259838, 1.020652192620172, 1.020370054495975, 1.019190961664661, 1.017113956130374, 1.015729718600066, 1.015130420748529, 1.014907826787519, 1.014955573423391, 1.015112020824417, 1.015208414784361, 1.015769493754813, 1.016121466131068, 1.016326774043065, 1.016396780186255, 1.016083701967358, 1.0150659506425, 1.014426118595012, 1.014765299503109, 1.016100765941605, 1.017650832406737)
    "F" = @(1.015869691325769, 1.017568849712941, 1.018667787984145, 1.019129662561586, 1.019207206585259, 1.018673960032807, 1.016444044597401, 1.012510224125855, 1.009884199313228, 1.008746164484176, 1.008323295421758, 1.008414009317634, 1.008711213157949, 1.008894309785883, 1.009959706119626, 1.01062774024205, 1.011017302361961, 1.011150117855384, 1.010556075975648, 1.008623698315916, 1.007407850025688, 1.008052481429773, 1.010588458252884, 1.0135277869713)
    "I" = @(1.026927994660535, 1.027077371256606, 1.027172427672561, 1.027212006320759, 1.027218629287371, 1.027172958029569, 1.026978808537782, 1.026624434729114, 1.026379940321444, 1.026272114562841, 1.026231768917422, 1.026240436518577, 1.026268785585074, 1.026286213369844, 1.026387055092431, 1.026449786084819, 1.026486187099213, 1.026498566866866, 1.026443075184663, 1.026260445618759, 1.026143916006149, 1.026205852059711, 1.026446108135329, 1.026717501823305)
    "J" = @(1.022720250132159, 1.023411121460803, 1.023857311615084, 1.02404468744639, 1.024076136824634, 1.023859816134697, 1.022953909851761, 1.02135103634721, 1.020277987655714, 1.01981227181354, 1.019639120807948, 1.019676269702652, 1.019797962438672, 1.019872919656319, 1.020308872841237, 1.020582045017059, 1.020741277741338, 1.020795554388035, 1.020552746993917, 1.019762131482711, 1.019264094038661, 1.019528203243208, 1.020565985830769, 1.021766200210451)
    "K" = @(1.026184222918567, 1.026758528330028, 1.027128946393864, 1.027284384282853, 1.027310466237882, 1.027131024485271, 1.026378560019735, 1.025043441239661, 1.024147158011155, 1.023757578049253, 1.023612647056569, 1.023643745379771, 1.023745602575131, 1.023808330495616, 1.024172981829881, 1.02440132027778, 1.024534363158259, 1.024579703097855, 1.024376836528303, 1.023715614349976, 1.023298583362464, 1.02351978233495, 1.024387900122936, 1.025389692579326)
    "L" = @(1.021732345631859, 1.022436097497875, 1.022890942459054, 1.023082033227602, 1.023114110837523, 1.022893496317679, 1.021970291956028, 1.020339399011131, 1.019249338982973, 1.018776655610959, 1.018600976672673, 1.018638665093039, 1.01876213604855, 1.018838196859618, 1.01928069499752, 1.019558079512635, 1.019719807654928, 1.019774941685838, 1.019528325568332, 1.018725779789709, 1.018220589261873, 1.018488457386481, 1.019541770297612, 1.020761513477644)
    "M" = @(1.018723559631315, 1.020225499997718, 1.021196389025617, 1.021604326263554, 1.021672807756949, 1.021201840776527, 1.019231351786962, 1.015751382053041, 1.013425777247569, 1.012417332953745, 1.012042526721776, 1.012122934254481, 1.012386355995229, 1.012548628738959, 1.013492673269974, 1.014084455419502, 1.014429493371646, 1.014547119149053, 1.014020977190884, 1.012308791140859, 1.011230965108877, 1.011802467783703, 1.014049660697454, 1.016651996916855)
}

foreach ($col in $newValues.Keys) {
    $vals = $newValues[$col]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $row = $i + 2
        $ws.Range("$col$row").Value = $vals[$i]
    }
}
